$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Feb 2018 updates: expanded product_collection line-item rows (A,B,C) plus
# a fixed-width numeric display format on the prod_id columns (B & C).
#
# Columns: row, prod_id_key (A), prod_id_value (B), next_month_prod_id (C), styleFlag
#   styleFlag 1 -> number format only
#   styleFlag 2 -> number format + yellow fill (matches legacy highlighted rows)
$data = @(
  @(2, "power_moves3", 69026938898, 138427301906, 2),
  @(3, "power_moves5", 69026316306, 138427203602, 1),
  @(4, "go_time3", 52386937778, 138427301906, 2),
  @(5, "go_time5", 78480408594, 138427203602, 1),
  @(6, "power_moves3_auto", 78541520914, 138427301906, 2),
  @(7, "power_moves5_auto", 78657093650, 138427203602, 1),
  @(8, "go_time3_auto", 91049066514, 138427301906, 2),
  @(9, "go_time5_auto1", 91049230354, 138427203602, 1),
  @(10, "go_time5_auto2", 91049197586, 138427203602, 1),
  @(11, "heart_soul_3item", 91236466706, 138427301906, 2),
  @(12, "heart_soul_3item_auto", 109303332882, 138427301906, 2),
  @(13, "heart_soul_5item", 109301366802, 138427203602, 1),
  @(14, "heart_soul_5item_auto", 109301366802, 138427203602, 1),
  @(15, "go_time_3item", 52386037778, 138427301906, 2),
  @(16, "fit_fierce_5item", 91235975186, 138427203602, 1),
  @(17, "modern_muse_5itemauto", 91236368402, 138427203602, 1),
  @(18, "power_move_3_item", 69026938898, 138427301906, 2),
  @(19, "fit_fierce_3item", 91236171794, 138427301906, 2),
  @(20, "modern_muse_5item", 91236368402, 138427203602, 1),
  @(21, "modern_muse_3item", 91236466706, 138427301906, 2),
  @(22, "fit_fierce_3itemauto", 91236171794, 138427301906, 2),
  @(23, "modern_muse_3itemauto", 91236466706, 138427301906, 2)
)

$numFmt = "###########00000"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $row[0]
    $key = $row[1]
    $val = $row[2]
    $nextVal = $row[3]
    $flag = $row[4]

    $ws.Cells.Item($r, 1).Value = $key
    $ws.Cells.Item($r, 2).Value = $val
    $ws.Cells.Item($r, 3).Value = $nextVal

    $ws.Cells.Item($r, 2).NumberFormat = $numFmt
    $ws.Cells.Item($r, 3).NumberFormat = $numFmt

    if ($flag -eq 2) {
        $ws.Cells.Item($r, 3).Interior.Color = 65535
    }
}

# Trailing blank data row (24) keeps the same numeric formatting as the
# rest of the table, no fill, no values.
$ws.Cells.Item(24, 2).NumberFormat = $numFmt
$ws.Cells.Item(24, 3).NumberFormat = $numFmt

# Column widths (characters, engine adds the standard padding offset)
$ws.Columns.Item(1).ColumnWidth = 22.33
$ws.Columns.Item(2).ColumnWidth = 12.5
$ws.Columns.Item(3).ColumnWidth = 18

# Match the saved selection/active cell from the edited workbook
$null = $ws.Range("C21:C23").Select()
